# Insert a new data row above the current row 158 (shifting the existing
# rows 158-281 down to 159-282) and populate the new row with its values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(158).Insert()

$ws.Range("A158").Value = 10
$ws.Range("B158").Value = "Vega Modelo de Temuco"
$ws.Range("C158").Value = "La Araucanía"
$ws.Range("D158").Value = 44669
$ws.Range("E158").Value = 9
$ws.Range("F158").Value = 100112017
$ws.Range("G158").Value = "Apio"
$ws.Range("H158").Value = "Americana (o)"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 85
$ws.Range("K158").Value = 9000
$ws.Range("L158").Value = 10000
$ws.Range("M158").Value = 9412
$ws.Range("N158").Value = "$/docena de matas"
$ws.Range("O158").Value = "Provincia del Elquí"
$ws.Range("P158").Value = 1569
$ws.Range("Q158").Value = 6
$ws.Range("R158").Value = "Hortaliza"
